$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B22").Value = 25
$ws.Range("B23").Value = 875000
$ws.Range("B24").Value = 3571428.571428572
$ws.Range("B31").Value = -5000000
$ws.Range("B34").Value = -403571.4285714282
$ws.Range("B35").Value = -403571.4285714282
